$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to be treated as text (avoids "+9989..." being
    # auto-converted to a number), then restore the original "s=2"
    # style by pasting formats from a known plain-text/number cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Update existing row 2 ---
$ws.Cells.Item(2,2).Value = "vcxvcvcxv"
Set-TextValue $ws.Cells.Item(2,3) "+998939950036"
$ws.Cells.Item(2,4).Value = "STORE_OWNER: fdsfdsf"
$ws.Cells.Item(2,5).Value = "ACCEPTED"
$ws.Cells.Item(2,6).Value = 50000
$ws.Cells.Item(2,7).Value = 3434
$ws.Cells.Item(2,8).Value = 44882.50939576389
$ws.Cells.Item(2,9).Value = 44882.51629447917

# --- Update existing row 3 ---
$ws.Cells.Item(3,2).Value = "recipient"
Set-TextValue $ws.Cells.Item(3,3) "+998939950202"
$ws.Cells.Item(3,4).Value = "STORE_OWNER: undefined"
$ws.Cells.Item(3,5).Value = "ACCEPTED"
$ws.Cells.Item(3,7).Value = 557577
$ws.Cells.Item(3,8).Value = 44882.511876793986
$ws.Cells.Item(3,9).Value = 44882.51630575232

# --- Update existing row 4 ---
$ws.Cells.Item(4,2).Value = "recipient"
Set-TextValue $ws.Cells.Item(4,3) "+998939950202"
$ws.Cells.Item(4,4).Value = "STORE_OWNER: undefined"
$ws.Cells.Item(4,5).Value = "ACCEPTED"
$ws.Cells.Item(4,7).Value = 557577
$ws.Cells.Item(4,8).Value = 44882.511935208335
$ws.Cells.Item(4,9).Value = 44882.51632061343

# --- Add new rows 5, 6, 7 by copying formatting from row 4 ---
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I7").PasteSpecial(-4122)

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "recipient"
Set-TextValue $ws.Cells.Item(5,3) "+998939950202"
$ws.Cells.Item(5,4).Value = "STORE_OWNER: undefined"
$ws.Cells.Item(5,5).Value = "ACCEPTED"
$ws.Cells.Item(5,6).Value = 50000
$ws.Cells.Item(5,7).Value = 557577
$ws.Cells.Item(5,8).Value = 44882.51194377315
$ws.Cells.Item(5,9).Value = 44882.51632913195

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "recipient"
Set-TextValue $ws.Cells.Item(6,3) "+998939950202"
$ws.Cells.Item(6,4).Value = "STORE_OWNER: undefined"
$ws.Cells.Item(6,5).Value = "ACCEPTED"
$ws.Cells.Item(6,6).Value = 50000
$ws.Cells.Item(6,7).Value = 557577
$ws.Cells.Item(6,8).Value = 44882.51195315972
$ws.Cells.Item(6,9).Value = 44882.51633701389

# --- Row 7 ---
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "ssss"
Set-TextValue $ws.Cells.Item(7,3) "+998939950036"
$ws.Cells.Item(7,4).Value = "STORE_OWNER: 232"
$ws.Cells.Item(7,5).Value = "ACCEPTED"
$ws.Cells.Item(7,6).Value = 50000
$ws.Cells.Item(7,7).Value = 433443
$ws.Cells.Item(7,8).Value = 44882.51458509259
$ws.Cells.Item(7,9).Value = 44882.51634814814

# --- Restore the original (non "@") number format / style on the
#     phone-number cells so their style index matches the other
#     cells in the row ("s=2") instead of a freshly created "@" style.
#     (PasteSpecial to a multi-area range only affects the first area,
#     so do it one cell at a time.)
foreach ($r in 2,3,4,5,6,7) {
    $ws.Cells.Item(4,1).Copy()
    $ws.Cells.Item($r,3).PasteSpecial(-4122)
}
